# Fix Training Data Issue (#48)
# The "Date" column (BF) values were off by one day because of how the
# NBA stats site displayed the date. Replace the mangled label
# "6-21-2011-12" with the correct ISO date "2012-06-21" on every data
# row (rows 2-31). The leading apostrophe forces Excel to keep the
# value as literal text instead of re-interpreting it as a date serial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $ws.Range("BF$row").Value = "'2012-06-21"
}
